$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp shown in F1.
$ws.Range("F1").Value = "Last status check on: 25.02.2022 09:30"

# Row 8 (Benzina Albert Modrice) values were recorded as plain text by the
# previous scraping script; re-write them as real numbers like every other
# row in the sheet.

# D8 (Delta Cena): text "+0.4" -> numeric 0.4
$ws.Range("D8").Value = 0.4

# E8 (Old Datum): text "2022-02-25 09:17:23" -> numeric Excel date/time,
# using the same date/time number format as the other rows in column E.
$ws.Range("E8").NumberFormat = $ws.Range("E7").NumberFormat
$ws.Range("E8").Value = 44617.38707175926
